$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 (shifts old rows 36-81 down to 37-82,
# growing the sheet dimension to A1:R82).
$ws.Rows.Item(36).Insert()

# Populate the freshly inserted row 36 with the new weekly price record.
# Columns A, B, C, E, F, G, H, I, N, O, Q, R keep the same values the
# (now shifted-down) row used to have; D, J, K, L, M, P carry the new data.
$ws.Range("A36").Value = 7
$ws.Range("B36").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C36").Value = "Ñuble"
$ws.Range("D36").Value = 44874
$ws.Range("E36").Value = 16
$ws.Range("F36").Value = 100112022
$ws.Range("G36").Value = "Arveja Verde"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 60
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 16000
$ws.Range("M36").Value = 15500
$ws.Range("N36").Value = "`$/saco 25 kilos"
$ws.Range("O36").Value = "Región del Maule"
$ws.Range("P36").Value = 620
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
